# Update TOTAL_SUBSTATION_LOAD / CONTESTABLE_ENERGY / ACTUAL_ENERGY figures
# with the latest DAP file numbers (rows 2-25). Rows 16-18 no longer have a
# TOTAL_SUBSTATION_LOAD (column B) reading, so those cells are cleared.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 33061
$ws.Range("C2").Value = 5005.48337366473
$ws.Range("D2").Value = 28055.51662633527
$ws.Range("B3").Value = 31828
$ws.Range("C3").Value = 4922.298456099087
$ws.Range("D3").Value = 26905.70154390091
$ws.Range("B4").Value = 30574
$ws.Range("C4").Value = 4850.625419333438
$ws.Range("D4").Value = 25723.37458066656
$ws.Range("B5").Value = 29643
$ws.Range("C5").Value = 4797.274898040722
$ws.Range("D5").Value = 24845.72510195928
$ws.Range("B6").Value = 29848
$ws.Range("C6").Value = 4859.095012476347
$ws.Range("D6").Value = 24988.90498752365
$ws.Range("B7").Value = 30796
$ws.Range("C7").Value = 4960.889853826959
$ws.Range("D7").Value = 25835.11014617304
$ws.Range("B8").Value = 31613
$ws.Range("C8").Value = 5211.449593293681
$ws.Range("D8").Value = 26401.55040670632
$ws.Range("B9").Value = 34939
$ws.Range("C9").Value = 6049.649550323782
$ws.Range("D9").Value = 28889.35044967622
$ws.Range("B10").Value = 41934
$ws.Range("C10").Value = 8417.127547097272
$ws.Range("D10").Value = 33516.87245290272
$ws.Range("B11").Value = 43152
$ws.Range("C11").Value = 13902.175
$ws.Range("D11").Value = 29249.825
$ws.Range("B12").Value = 44379
$ws.Range("C12").Value = 15210.828
$ws.Range("D12").Value = 29168.172
$ws.Range("B13").Value = 45163
$ws.Range("C13").Value = 15097.205
$ws.Range("D13").Value = 30065.795
$ws.Range("B14").Value = 45347
$ws.Range("C14").Value = 14468.754
$ws.Range("D14").Value = 30878.246
$ws.Range("B15").Value = 47292
$ws.Range("C15").Value = 14615.9345
$ws.Range("D15").Value = 32676.0655
$ws.Range("B16").ClearContents()
$ws.Range("C16").Value = 15642.9465
$ws.Range("D16").Value = 49763.3545
$ws.Range("B17").ClearContents()
$ws.Range("C17").Value = 15864.0115
$ws.Range("D17").Value = 28812.386
$ws.Range("B18").ClearContents()
$ws.Range("C18").Value = 16275.092
$ws.Range("D18").Value = 96.81750000000102
$ws.Range("C19").Value = 16265.2175
$ws.Range("D19").Value = 39727.5055
$ws.Range("C20").Value = 15390.046
$ws.Range("D20").Value = 39165.7675
$ws.Range("C21").Value = 13874.261
$ws.Range("D21").Value = 37817.495
$ws.Range("C22").Value = 11951.5345
$ws.Range("D22").Value = 38010.227
$ws.Range("C23").Value = 9041.0965
$ws.Range("D23").Value = 38285.2065
$ws.Range("C24").Value = 6150.5125
$ws.Range("D24").Value = 0
$ws.Range("C25").Value = 5731.305
$ws.Range("D25").Value = 0